# Applies the "listo la base de datos y scripts sql" commit:
#  - extends the "Partidos" mini-table with id_persona_var / hr_partido+Hora,
#    and renames lf_partido -> dt_partido
#  - that pushes the "Penaltis" and "TipoPersonaTecnica" mini-tables (and the
#    comment anchored on the latter) down by two rows
#  - recolors the table-header highlight from green to yellow
#  - leaves the cursor on K10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Move the "Penaltis" and "TipoPersonaTecnica" blocks (cols I:J only)
#        down two rows to make room for the Partidos table growing.       ---
$ws.Range("I14:J14").UnMerge()
$ws.Range("I20:J20").UnMerge()

# lower block: rows 20-22 -> 22-24 (bottom-up, ranges overlap)
$ws.Range("I24").Value = $ws.Range("I22").Value2
$ws.Range("J24").Value = $ws.Range("J22").Value2
$ws.Range("I23").Value = $ws.Range("I21").Value2
$ws.Range("J23").Value = $ws.Range("J21").Value2
$ws.Range("I22").Value = $ws.Range("I20").Value2
$ws.Range("J22").Value = $ws.Range("J20").Value2

# middle block: rows 14-17 -> 16-19 (bottom-up, ranges overlap)
$ws.Range("I19").Value = $ws.Range("I17").Value2
$ws.Range("J19").Value = $ws.Range("J17").Value2
$ws.Range("I18").Value = $ws.Range("I16").Value2
$ws.Range("J18").Value = $ws.Range("J16").Value2
$ws.Range("I17").Value = $ws.Range("I15").Value2
$ws.Range("J17").Value = $ws.Range("J15").Value2
$ws.Range("I16").Value = $ws.Range("I14").Value2
$ws.Range("J16").Value = $ws.Range("J14").Value2

# the rows the content moved out of are now stale; rows 12/13 get new data below
$ws.Range("I14:J15").ClearContents()
$ws.Range("I20:J21").ClearContents()

# restore the merges at their new homes
$ws.Range("I16:J16").Merge()
$ws.Range("I22:J22").Merge()

# --- 3) Grow the "Partidos" table: rename lf_partido -> dt_partido, and add
#        id_persona_var (row 12) + hr_partido/Hora (row 13)                ---
$ws.Range("I6").Value = "dt_partido"
$ws.Range("I12").Value = "id_persona_var"
$ws.Range("J12").Value = "number"
$ws.Range("I13").Value = "hr_partido"
$ws.Range("J13").Value = "Hora"

# --- 4) Recolor the table-header highlight fill from green to yellow ---
$headerAreas = @("A1:B1","C1:D1","E1:F1","G1:H1","I1:J1","K1:L1","A5:B5","K5:L5", `
                 "E8:F8","A12:B12","C12:D12","E13:F13","I16:J16","G16:H16", `
                 "E17:F17","C20:D20","I22:J22")
foreach ($area in $headerAreas) {
    $ws.Range($area).Interior.Color = 65535
}

# --- 5) Move the comment that documents the TipoPersonaTecnica table ---
$oldComment = $ws.Range("I20").Comment
if ($oldComment -ne $null) {
    $commentText = $oldComment.Text()
    $oldComment.Delete()
    $ws.Range("I22").AddComment($commentText)
}

# --- 6) Selection left by the author when they saved ---
$ws.Range("K10").Select() | Out-Null

Write-Output "edit complete"
